$d = $word.ActiveDocument
$xmlns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Change 1: "...kept ideally." -> "...kept id" + "le" + "." (3 runs) ---
$orig1 = "The battery connected to it must be checked regularly since it may discharge itself when kept ideally."
$full = $d.Content.Text
$start1 = $full.IndexOf($orig1)
$end1 = $start1 + $orig1.Length
$rng1 = $d.Range($start1, $end1)

$xml1 = '<w:p ' + $xmlns + '>' `
    + '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>' `
    + '<w:t>The battery connected to it must be checked regularly since it may discharge itself when kept id</w:t></w:r>' `
    + '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>le</w:t></w:r>' `
    + '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>.</w:t></w:r>' `
    + '</w:p>'
$rng1.InsertXML($xml1)

# --- Change 2: split "On completion of the project we can say..." around
#     "project", wrapping it with gramStart/gramEnd proofErr markers ---
$orig2 = "On completion of the project we can say that this dispenser unit can be extremely helpful to old age people and to them who are alone at home the whole day and tend to forget to take medicine on time."
$full2 = $d.Content.Text
$start2 = $full2.IndexOf($orig2)
$end2 = $start2 + $orig2.Length
$rng2 = $d.Range($start2, $end2)

$xml2 = '<w:p ' + $xmlns + '>' `
    + '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>' `
    + '<w:t xml:space="preserve">On completion of the </w:t></w:r>' `
    + '<w:proofErr w:type="gramStart"/>' `
    + '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>project</w:t></w:r>' `
    + '<w:proofErr w:type="gramEnd"/>' `
    + '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>' `
    + '<w:t xml:space="preserve"> we can say that this dispenser unit can be extremely helpful to old age people and to them who are alone at home the whole day and tend to forget to take medicine on time.</w:t></w:r>' `
    + '</w:p>'
$rng2.InsertXML($xml2)
